# EleccionesPresidenciales2020 - rework header/data layout:
#  - Move the Presidentes/#Votos table from D5:E8 up to A1:B4
#  - Rename the "#Votos" header to "Voto"
#  - Turn the raw vote-share fractions into "* 100" formulas (percentage -> plain number)
#  - Number format for those cells becomes a plain 2-decimal number instead of a percentage
#  - Widen column A to fit the longer labels, move selection to B5, and set the page to portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move (cut+paste) the whole D5:E8 block to A1:B4 - this preserves the existing
# cell style (percentage number format) that travels with column E onto column B.
$ws.Range("D5:E8").Cut($ws.Range("A1:B4"))

# Remove whatever is left behind in the old D5:E8 block (Cut leaves empty,
# still-styled cells behind instead of truly clearing the block).
$ws.Range("D5:E8").Clear()

# Update the header text: "#Votos" -> "Voto"
$ws.Range("B1").Value = "Voto"

# Replace the raw fraction values with "<fraction>*100" formulas
$ws.Range("B2").Formula = "=0.433*100"
$ws.Range("B3").Formula = "=0.3028*100"
$ws.Range("B4").Formula = "=0.153*100"

# The values are no longer a fraction-of-one, so switch the number format from
# percentage to a plain number with 2 decimals.
$ws.Range("B2:B4").NumberFormat = "0.00"

# Column A needs to be wide enough for the candidate names.
$ws.Columns.Item(1).ColumnWidth = 19

# Match the saved selection/active cell.
$ws.Range("B5").Select() | Out-Null

# Page is now set to portrait orientation.
$ws.PageSetup.Orientation = 1
